$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 308; existing rows 308-386 shift down to 311-389
$ws.Range("A308:A310").EntireRow.Insert()

# Row 308
$ws.Cells.Item(308, 1).Value = 11
$ws.Cells.Item(308, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(308, 3).Value = "Bíobío"
$ws.Cells.Item(308, 4).Value = 44855
$ws.Cells.Item(308, 5).Value = 8
$ws.Cells.Item(308, 6).Value = "Fruta"
$ws.Cells.Item(308, 7).Value = 100101
$ws.Cells.Item(308, 8).Value = "Berries"
$ws.Cells.Item(308, 9).Value = 100112025
$ws.Cells.Item(308, 10).Value = "Frutilla"
$ws.Cells.Item(308, 11).Value = "Sin especificar"
$ws.Cells.Item(308, 12).Value = "Especial"
$ws.Cells.Item(308, 13).Value = 220
$ws.Cells.Item(308, 14).Value = 12000
$ws.Cells.Item(308, 15).Value = 13000
$ws.Cells.Item(308, 16).Value = 12455
$ws.Cells.Item(308, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(308, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(308, 19).Value = 1779
$ws.Cells.Item(308, 20).Value = 7

# Row 309
$ws.Cells.Item(309, 1).Value = 11
$ws.Cells.Item(309, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(309, 3).Value = "Bíobío"
$ws.Cells.Item(309, 4).Value = 44855
$ws.Cells.Item(309, 5).Value = 8
$ws.Cells.Item(309, 6).Value = "Fruta"
$ws.Cells.Item(309, 7).Value = 100101
$ws.Cells.Item(309, 8).Value = "Berries"
$ws.Cells.Item(309, 9).Value = 100112025
$ws.Cells.Item(309, 10).Value = "Frutilla"
$ws.Cells.Item(309, 11).Value = "Sin especificar"
$ws.Cells.Item(309, 12).Value = "Primera"
$ws.Cells.Item(309, 13).Value = 270
$ws.Cells.Item(309, 14).Value = 9000
$ws.Cells.Item(309, 15).Value = 10000
$ws.Cells.Item(309, 16).Value = 9556
$ws.Cells.Item(309, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(309, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(309, 19).Value = 1365
$ws.Cells.Item(309, 20).Value = 7

# Row 310
$ws.Cells.Item(310, 1).Value = 11
$ws.Cells.Item(310, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(310, 3).Value = "Bíobío"
$ws.Cells.Item(310, 4).Value = 44855
$ws.Cells.Item(310, 5).Value = 8
$ws.Cells.Item(310, 6).Value = "Fruta"
$ws.Cells.Item(310, 7).Value = 100101
$ws.Cells.Item(310, 8).Value = "Berries"
$ws.Cells.Item(310, 9).Value = 100112025
$ws.Cells.Item(310, 10).Value = "Frutilla"
$ws.Cells.Item(310, 11).Value = "Sin especificar"
$ws.Cells.Item(310, 12).Value = "Segunda"
$ws.Cells.Item(310, 13).Value = 170
$ws.Cells.Item(310, 14).Value = 7000
$ws.Cells.Item(310, 15).Value = 8000
$ws.Cells.Item(310, 16).Value = 7529
$ws.Cells.Item(310, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(310, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(310, 19).Value = 1076
$ws.Cells.Item(310, 20).Value = 7

